# Build site at 2023-04-12 14:53:07 UTC
# Update LOM3099 worksheet: fill in real content for Objetivos, Docentes
# responsaveis (3 separate rows), Programa resumido, Programa, Metodo,
# Criterio, Norma de recuperacao and Bibliografia.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert three new rows right after row 12 ("Docentes responsaveis:")
#    so that each professor gets their own row (13, 14, 15), pushing the
#    remaining rows down by three.
# ---------------------------------------------------------------------
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(13).Insert()

# The inserted rows copy column A's bold style; clear that since column A
# stays empty on these rows.
$ws.Range("A13:A15").Clear()

# Give the new B/C cells the same look (style) as other content cells.
$ws.Range("B2:C2").Copy() | Out-Null
$ws.Range("B13:C15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 2. Fill in the actual content (store each text in a variable first so
#    we never need to read a Value back from a Range).
# ---------------------------------------------------------------------

# Objetivos: real objective text instead of the placeholder professor name.
$objetivos = "Proporcionar ao aluno um conhecimento básico da mecânica dos corpos rígidos com ênfase na estática. Fornecer conhecimentos necessários para cálculo de reações de apoios e de esforços internos em estruturas isostáticas."
$ws.Range("B10").Value = $objetivos
$ws.Range("C10").Value = $objetivos

# Docentes responsáveis: one professor per row.
$prof1 = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("B13").Value = $prof1
$ws.Range("C13").Value = $prof1

$prof2 = "5840793 - Sérgio Schneider"
$ws.Range("B14").Value = $prof2
$ws.Range("C14").Value = $prof2

$prof3 = "7797767 - Viktor Pastoukhov"
$ws.Range("B15").Value = $prof3
$ws.Range("C15").Value = $prof3

# Programa resumido: (row 16 after the insert)
$resumido = "Estática de Partículas. Estática de Corpos Rígidos. Equilíbrio de Corpos Rígidos. Análise de Estruturas."
$ws.Range("B16").Value = $resumido
$ws.Range("C16").Value = $resumido

# Short syllabus: (row 17) no longer needs the extra row height.
$ws.Rows.Item(17).AutoFit()

# Programa: (row 18) real syllabus text; row height shrinks from 120 to 60.
$programa = "Mecânica e suas áreas: Corpos rígidos e corpos deformáveis (sólidos). Terminologia e metodologia básica. Estática de Partículas: Vetores, resultante de várias forças concorrentes, equilíbrio de uma partícula. Estática de Corpos Rígidos: Conceito de corpo rígido. Momento de uma força com relação a um ponto, sistemas equivalentes de forças, momento e binário. Apoios e vínculos. Diagrama de corpo livre. Reações de apoios e conexões para uma estrutura 2D. Equilíbrio de um corpo rígido em 2D. Reações estaticamente indeterminadas e vínculos parciais. Equilíbrio de um corpo rígido em 3D. Análise de Estruturas: análise do equilíbrio de estruturas, ação de múltiplas forças, forças internas, terceira Lei de Newton. Treliças: método dos nós, método das seções. Estruturas e Máquinas: transmissão e modificação de forças. Esforços internos em pórticos, vigas, cabos e eixos de transmissão."
$ws.Range("B18").Value = $programa
$ws.Range("C18").Value = $programa
$ws.Rows.Item(18).RowHeight = 60

# Método: (row 21) real evaluation method text.
$metodo = "Os alunos serão avaliados continuamente a qual serão considerados: provas escritas, exercícios propostos e seminários"
$ws.Range("B21").Value = $metodo
$ws.Range("C21").Value = $metodo

# Critério: (row 22) real grading criteria text.
$criterio = "Para compor a Nota Final (NF) serão consideradas as avaliações de Provas Escritas (P1 e P2) e Exercícios Propostos e Seminários (T) em que:NF = (P1 + P2 + T)/3.  Serão considerados aprovados os alunos que obtiverem: NF maior ou igual a 5,0. Serão considerados reprovados os alunos que obtiverem: NS menor que 3,0 Para os alunos em que NS é maior ou igual a 3,0 e menor que 5,0 será dada uma prova de recuperação (R)."
$ws.Range("B22").Value = $criterio
$ws.Range("C22").Value = $criterio

# Norma de recuperação: (row 23) real recovery-exam rule text.
$norma = ": A prova de Recuperação (R) irá compor a nota final de recuperação (NR) da seguinte forma: NR = (R + NF)/2. Serão considerados aprovados os alunos que obtiverem NF maior ou igual a 5,0."
$ws.Range("B23").Value = $norma
$ws.Range("C23").Value = $norma

# Bibliografia: (row 24) real bibliography text.
$biblio = "1. F.P. BEER, E.R. JOHNSTON, J.T. DeWOLF, D.MAZUREK. Estática e Mecânica dos Materiais. São Paulo: McGraw Hill, 2013, 728p.2. F.P. BEER, E.R. JOHNSTON, E. RUSSEL. Mecânica vetorial para engenheiros: Estática. São Paulo: McGraw Hill. 9a Ed., 2012, 626p. 3. HIBBELER, R.C. Mecânica para engenharia vol.1: estática. São Paulo: Pearson Prentice Hall, 12a Ed., 2011.4. MERIAM, J.L. KRAIGE, L.G. Mecânica para engenharia – Estática. Grupo GEN Editora LTC, 6a Ed., 2009, 364p. 5. RUIZ, C.C.de La P. Fundamentos de mecânica para engenharia – Estática. Grupo GEN Editora LTC, 2017, 306p."
$ws.Range("B24").Value = $biblio
$ws.Range("C24").Value = $biblio

# ---------------------------------------------------------------------
# 3. Column 1 used to share its width definition with column 2
#    ("min=1,max=2"); re-assert column 2's own width so it is tracked as
#    its own distinct column and column 1 no longer spans into it.
# ---------------------------------------------------------------------
$colBWidth = $ws.Columns.Item(2).ColumnWidth
$ws.Columns.Item(2).ColumnWidth = $colBWidth

